$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B (rows 28-72): values change from 6 to 0 ---
# Rows 29-31 also pick up the same cell formatting already used by the
# surrounding rows (32-72), which is applied first by copying the format
# from B32 before the values are overwritten.
$ws.Range("B32").Copy()
$ws.Range("B29:B31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 28; $r -le 72; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}

# --- Update the sheet view: scroll back to top, zoom to 130% ---
$ws.Range("A2").Select()
$excel.ActiveWindow.Zoom = 130
